# Reproduce manual station-row edits on the "Station" sheet (rows 111-171, columns E:F)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Station")

$ws.Cells.Item(111, 5).Value = 372
$ws.Cells.Item(111, 6).Value = 59
$ws.Cells.Item(112, 5).Value = 356
$ws.Cells.Item(112, 6).Value = 59
$ws.Cells.Item(113, 5).Value = 344
$ws.Cells.Item(113, 6).Value = 71
$ws.Cells.Item(114, 5).Value = 332
$ws.Cells.Item(114, 6).Value = 83
$ws.Cells.Item(115, 5).Value = 320
$ws.Cells.Item(115, 6).Value = 95
$ws.Cells.Item(116, 5).Value = 308
$ws.Cells.Item(116, 6).Value = 107
$ws.Cells.Item(117, 5).Value = 296
$ws.Cells.Item(117, 6).Value = 119
$ws.Cells.Item(118, 5).Value = 284
$ws.Cells.Item(118, 6).Value = 131
$ws.Cells.Item(119, 5).Value = 284
$ws.Cells.Item(119, 6).Value = 144
$ws.Cells.Item(120, 5).Value = 284
$ws.Cells.Item(120, 6).Value = 158
$ws.Cells.Item(121, 5).Value = 284
$ws.Cells.Item(121, 6).Value = 170
$ws.Cells.Item(122, 5).Value = 284
$ws.Cells.Item(122, 6).Value = 182
$ws.Cells.Item(123, 5).Value = 284
$ws.Cells.Item(123, 6).Value = 194
$ws.Cells.Item(124, 5).Value = 284
$ws.Cells.Item(124, 6).Value = 206
$ws.Cells.Item(125, 5).Value = 284
$ws.Cells.Item(125, 6).Value = 218
$ws.Cells.Item(126, 5).Value = 284
$ws.Cells.Item(126, 6).Value = 231
$ws.Cells.Item(127, 5).Value = 284
$ws.Cells.Item(127, 6).Value = 246
$ws.Cells.Item(128, 5).Value = 284
$ws.Cells.Item(128, 6).Value = 260
$ws.Cells.Item(129, 5).Value = 284
$ws.Cells.Item(129, 6).Value = 274
$ws.Cells.Item(130, 5).Value = 284
$ws.Cells.Item(130, 6).Value = 288
$ws.Cells.Item(131, 5).Value = 284
$ws.Cells.Item(131, 6).Value = 302
$ws.Cells.Item(132, 5).Value = 284
$ws.Cells.Item(132, 6).Value = 318
$ws.Cells.Item(133, 5).Value = 284
$ws.Cells.Item(133, 6).Value = 334
$ws.Cells.Item(134, 5).Value = 284
$ws.Cells.Item(134, 6).Value = 348
$ws.Cells.Item(135, 5).Value = 284
$ws.Cells.Item(135, 6).Value = 362
$ws.Cells.Item(136, 5).Value = 300
$ws.Cells.Item(136, 6).Value = 362
$ws.Cells.Item(137, 5).Value = 316
$ws.Cells.Item(137, 6).Value = 362
$ws.Cells.Item(138, 5).Value = 332
$ws.Cells.Item(138, 6).Value = 362
$ws.Cells.Item(139, 5).Value = 347
$ws.Cells.Item(139, 6).Value = 362
$ws.Cells.Item(140, 5).Value = 363
$ws.Cells.Item(140, 6).Value = 362
$ws.Cells.Item(141, 5).Value = 370
$ws.Cells.Item(141, 6).Value = 374
$ws.Cells.Item(142, 5).Value = 377
$ws.Cells.Item(142, 6).Value = 386
$ws.Cells.Item(143, 5).Value = 384
$ws.Cells.Item(143, 6).Value = 398
$ws.Cells.Item(144, 5).Value = 391
$ws.Cells.Item(144, 6).Value = 410
$ws.Cells.Item(145, 5).Value = 398
$ws.Cells.Item(145, 6).Value = 422
$ws.Cells.Item(146, 5).Value = 405
$ws.Cells.Item(146, 6).Value = 434
$ws.Cells.Item(147, 5).Value = 405
$ws.Cells.Item(147, 6).Value = 448
$ws.Cells.Item(148, 5).Value = 405
$ws.Cells.Item(148, 6).Value = 462
$ws.Cells.Item(149, 5).Value = 405
$ws.Cells.Item(149, 6).Value = 476
$ws.Cells.Item(150, 5).Value = 405
$ws.Cells.Item(150, 6).Value = 490
$ws.Cells.Item(151, 5).Value = 405
$ws.Cells.Item(151, 6).Value = 504
$ws.Cells.Item(152, 5).Value = 405
$ws.Cells.Item(152, 6).Value = 518
$ws.Cells.Item(153, 5).Value = 405
$ws.Cells.Item(153, 6).Value = 532
$ws.Cells.Item(154, 5).Value = 405
$ws.Cells.Item(154, 6).Value = 546
$ws.Cells.Item(155, 5).Value = 405
$ws.Cells.Item(155, 6).Value = 560
$ws.Cells.Item(156, 5).Value = 405
$ws.Cells.Item(156, 6).Value = 574
$ws.Cells.Item(157, 5).Value = 405
$ws.Cells.Item(157, 6).Value = 588
$ws.Cells.Item(158, 5).Value = 421
$ws.Cells.Item(158, 6).Value = 604
$ws.Cells.Item(159, 5).Value = 272
$ws.Cells.Item(159, 6).Value = 362
$ws.Cells.Item(160, 5).Value = 284
$ws.Cells.Item(160, 6).Value = 373
$ws.Cells.Item(161, 5).Value = 284
$ws.Cells.Item(161, 6).Value = 384
$ws.Cells.Item(162, 5).Value = 284
$ws.Cells.Item(162, 6).Value = 406
$ws.Cells.Item(163, 5).Value = 272
$ws.Cells.Item(163, 6).Value = 416
$ws.Cells.Item(164, 5).Value = 260
$ws.Cells.Item(164, 6).Value = 426
$ws.Cells.Item(165, 5).Value = 248
$ws.Cells.Item(165, 6).Value = 436
$ws.Cells.Item(166, 5).Value = 232
$ws.Cells.Item(166, 6).Value = 446
$ws.Cells.Item(167, 5).Value = 216
$ws.Cells.Item(167, 6).Value = 446
$ws.Cells.Item(168, 5).Value = 200
$ws.Cells.Item(168, 6).Value = 446
$ws.Cells.Item(169, 5).Value = 183
$ws.Cells.Item(169, 6).Value = 446
$ws.Cells.Item(170, 5).Value = 167
$ws.Cells.Item(170, 6).Value = 429
$ws.Cells.Item(171, 5).Value = 149
$ws.Cells.Item(171, 6).Value = 411

# Restore the view/selection state recorded after the edit
$ws.Activate()
$ws.Range("E160:F171").Select()
